# Updates crypto price/volume figures (and restores the WrappedeETH /
# Binance-PegBSC-USD row order) per the Nov 3 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=1 (#), B=2 (Coin), C=3 (Link), D=4 (Price), E=5 (Volume(1h))
# D-column values that parse as plain numbers need NumberFormat "@" forced
# first, so Excel keeps them as text (matching the original inline-string
# cells) instead of silently coercing them to numeric values.

$ws.Cells.Item(2, 4).Value = '68.447.85'
$ws.Cells.Item(2, 5).Value = '  -1.88%  '
$ws.Cells.Item(3, 4).Value = '2.451.45'
$ws.Cells.Item(3, 5).Value = '  -2.27%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '564.28'
$ws.Cells.Item(5, 5).Value = '  -2.03%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '163.84'
$ws.Cells.Item(6, 5).Value = '  -1.96%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$ws.Cells.Item(8, 5).Value = '  -1.38%  '
$ws.Cells.Item(9, 5).Value = '  -6.62%  '
$ws.Cells.Item(10, 5).Value = '  -2.03%  '
$ws.Cells.Item(11, 5).Value = '  -4.07%  '
$ws.Cells.Item(12, 5).Value = '  -2.67%  '
$ws.Cells.Item(13, 4).Value = '2.902.50'
$ws.Cells.Item(13, 5).Value = '  -1.84%  '
$ws.Cells.Item(14, 4).Value = '68.406.77'
$ws.Cells.Item(14, 5).Value = '  -1.81%  '
$ws.Cells.Item(15, 5).Value = '  -4.40%  '
$ws.Cells.Item(16, 5).Value = '  -5.19%  '
$ws.Cells.Item(17, 4).Value = '2.482.41'
$ws.Cells.Item(17, 5).Value = '  -0.09%  '
$ws.Cells.Item(18, 5).Value = '  -2.24%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '344.36'
$ws.Cells.Item(19, 5).Value = '  -1.53%  '
$ws.Cells.Item(20, 5).Value = '  -4.66%  '
$ws.Cells.Item(21, 5).Value = '  -2.09%  '
$ws.Cells.Item(22, 5).Value = '  -3.32%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.999'
$ws.Cells.Item(23, 5).Value = '  -0.07%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '68.12'
$ws.Cells.Item(24, 5).Value = '  -3.46%  '
$ws.Cells.Item(25, 5).Value = '  -5.06%  '
$ws.Cells.Item(26, 2).Value = 'WrappedeETH'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(26, 4).Value = '2.581.51'
$ws.Cells.Item(26, 5).Value = '  -0.38%  '
$ws.Cells.Item(27, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '1.03'
$ws.Cells.Item(27, 5).Value = '  +3.61%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '8.24'
$ws.Cells.Item(28, 5).Value = '  -6.62%  '
$ws.Cells.Item(30, 5).Value = '  -7.06%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '437.92'
$ws.Cells.Item(31, 5).Value = '  -4.86%  '
$ws.Cells.Item(32, 5).Value = '  -3.42%  '
$ws.Cells.Item(33, 5).Value = '  +0.02%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.68'
$ws.Cells.Item(34, 5).Value = '  -2.76%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '3.07'
$ws.Cells.Item(35, 5).Value = '  +103.96%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '156.72'
$ws.Cells.Item(36, 5).Value = '  -0.94%  '
$ws.Cells.Item(37, 5).Value = '  -0.38%  '
$ws.Cells.Item(38, 5).Value = '  +0.00%  '
$ws.Cells.Item(39, 5).Value = '  -6.21%  '
$ws.Cells.Item(41, 5).Value = '  -3.72%  '
$ws.Cells.Item(42, 5).Value = '  -4.28%  '
$ws.Cells.Item(43, 5).Value = '  -4.45%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.14'
$ws.Cells.Item(44, 5).Value = '  +5.31%  '
$ws.Cells.Item(45, 5).Value = '  -5.40%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '135.32'
$ws.Cells.Item(46, 5).Value = '  -4.41%  '
$ws.Cells.Item(47, 5).Value = '  -3.13%  '
$ws.Cells.Item(48, 5).Value = '  -2.39%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.487'
$ws.Cells.Item(49, 5).Value = '  -6.42%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.565'
$ws.Cells.Item(50, 5).Value = '  -2.52%  '
$ws.Cells.Item(51, 5).Value = '  -1.46%  '
